# edit.ps1 -- applies the BlurAnalysis.docx text revisions described in the
# commit "small changes to blurAnalysis".
#
# The substantive, reproducible changes are:
#   1. Summary paragraph: reword "The size of this area is decided by the
#      user," -> "The user decides the size of this area,"
#   2. Introduction paragraph: "allows the client to blur the picture." ->
#      "allows the client to apply a blur effect to a bitmap picture."
#   3. Run-time paragraph: rewrite the "...pixels large, and must do this
#      for each pixel in the picture.  Therefore, the number..." sentence
#      to "...pixels large, where n is the threshold or range chosen by
#      the user. This long iteration is done to each pixel in the
#      picture,  therefore, the number..."
#   4. Same paragraph: "which is very fast" -> "which is very slow"
#   5. Add a "_GoBack" bookmark (empty range) to the blank paragraph just
#      before the document's trailing blank paragraph.
#
# (The diff's many <w:proofErr w:type="gramStart|gramEnd|spellStart|
# spellEnd"/> markers are artifacts Word's background proofing tool stamps
# into the markup; they carry no visible text and aren't reachable through
# the Word object model, so they are not something this script produces.)

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Host ("NOT FOUND: " + $oldText)
    }
    return $ok
}

# 1. Summary paragraph.
Replace-Text "The size of this area is decided by the user," "The user decides the size of this area,"

# 2. Introduction paragraph.
Replace-Text "allows the client to blur the picture." "allows the client to apply a blur effect to a bitmap picture."

# 3. Run-time paragraph -- rewrite the "and must do this ... Therefore," sentence.
Replace-Text "pixels large, and must do this for each pixel in the picture.  Therefore, the number" "pixels large, where n is the threshold or range chosen by the user. This long iteration is done to each pixel in the picture,  therefore, the number"

# 4. "very fast" -> "very slow" (run-time growth paragraph).
Replace-Text "which is very fast" "which is very slow"

# 5. Add the "_GoBack" bookmark to the blank paragraph before the trailing
#    blank paragraph (Word stamps this on the last edit location).
$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count - 1)
$d.Bookmarks.Add("_GoBack", $target.Range)

Write-Host "done"
